$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking values must be forced to remain as text
# (otherwise Excel would auto-convert them to numbers and drop formatting,
# e.g. trailing zeros in "143.70").
$textCells = @("D5", "D6", "D10", "D14", "D15", "D20", "D25", "D29", "D41", "D46", "D47", "D48", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "25.923.93"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.623.62"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "213.51"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("D10").Value = "18.24"
$ws.Range("E10").Value = "  -6.77%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").Value = "1.848.91"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.646.15"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.18"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("D16").Value = "25.910.83"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -3.49%  "
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -3.93%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "191.17"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  -2.89%  "
$ws.Range("E22").Value = "  -3.66%  "
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "143.70"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("D29").Value = "15.15"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("E33").Value = "  -5.47%  "
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("D36").Value = "1.118.97"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("E37").Value = "  -6.57%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").Value = "97.88"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("D43").Value = "1.759.23"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  -5.74%  "
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").Value = "0.0529"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "54.38"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("D48").Value = "1.47"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").Value = "7.45"
$ws.Range("E51").Value = "  -3.94%  "
